$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (position + assignment) for Anakin Skywalker / Darth Vader
$ws.Range("A13").Value = "askywalker"
$ws.Range("B13").Value = "Anakin"
$ws.Range("C13").Value = "Skywalker"
$ws.Range("D13").Value = "vader@sith.com"
$ws.Range("E13").Value = 2233355
$ws.Range("F13").Value = "TINO-NS"
$ws.Range("G13").Value = "200-1234"

# Turn the email cell into a hyperlink, matching the pattern used by the other rows
$ws.Hyperlinks.Add($ws.Range("D13"), "../../../Documents/vader@sith.com")

# Re-apply the shared "Hyperlink" cell style so D13 matches the other linked e-mail cells
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122)

# Update the current selection to cover the newly added row
$ws.Range("E12:G13").Select()
